$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDataBatch16")

# Update the "username" column values (column E) for rows 2-4
$ws.Range("E2").Value = "fire963"
$ws.Range("E3").Value = "water159"
$ws.Range("E4").Value = "soil369"
